$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 522.2759
$ws.Range("I33").Value = 166.38095
$ws.Range("J33").Value = 1456.5
$ws.Range("K33").Value = 166.38095
$ws.Range("L33").Value = 1456.5
$ws.Range("M33").Value = 62.61904999999999
$ws.Range("N33").Value = -1914.5

$ws.Range("H53").Value = 180.5
$ws.Range("I53").Value = 57.6
$ws.Range("J53").Value = 385.33334
$ws.Range("K53").Value = 57.6
$ws.Range("L53").Value = 385.33334
$ws.Range("M53").Value = 579.4
$ws.Range("N53").Value = -1659.33334

$ws.Range("H87").Value = 43569.332
$ws.Range("J87").Value = 50354
$ws.Range("L87").Value = 50354
$ws.Range("N87").Value = -52850

$ws.Range("H90").Value = 43569.332
$ws.Range("J90").Value = 50354
$ws.Range("L90").Value = 151062
$ws.Range("N90").Value = -163542

$ws.Range("H100").Value = 3976.9565
$ws.Range("I100").Value = 2673
$ws.Range("J100").Value = 6421.875
$ws.Range("K100").Value = 2673
$ws.Range("L100").Value = 6421.875
$ws.Range("M100").Value = -2132
$ws.Range("N100").Value = -7503.875

$ws.Range("H132").Value = 206766.33
$ws.Range("I132").Value = 2571.5386
$ws.Range("J132").Value = 1003126
$ws.Range("K132").Value = 7714.6158
$ws.Range("L132").Value = 3009378
$ws.Range("M132").Value = -5184.6158
$ws.Range("N132").Value = -3014438

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 5502.5
$ws.Range("I63").Value = 5005
$ws.Range("K63").Value = 5005
$ws.Range("M63").Value = -4319

$ws.Range("H66").Value = 5502.5
$ws.Range("I66").Value = 5005
$ws.Range("K66").Value = 25025
$ws.Range("M66").Value = -21593

$ws.Range("H80").Value = 12833.25
$ws.Range("J80").Value = 12833.25
$ws.Range("L80").Value = 12833.25
$ws.Range("N80").Value = -14829.25

$ws.Range("H83").Value = 12833.25
$ws.Range("J83").Value = 12833.25
$ws.Range("L83").Value = 38499.75
$ws.Range("N83").Value = -48483.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 5133.8096
$ws.Range("I80").Value = 11329.777
$ws.Range("J80").Value = 486.83334
$ws.Range("K80").Value = 11329.777
$ws.Range("L80").Value = 486.83334
$ws.Range("M80").Value = -10331.777
$ws.Range("N80").Value = -2482.83334

$ws.Range("H83").Value = 5133.8096
$ws.Range("I83").Value = 11329.777
$ws.Range("J83").Value = 486.83334
$ws.Range("K83").Value = 56648.885
$ws.Range("L83").Value = 2434.1667
$ws.Range("M83").Value = -51656.885
$ws.Range("N83").Value = -12418.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 46255.5
$ws.Range("J14").Value = 46255.5
$ws.Range("L14").Value = 46255.5
$ws.Range("N14").Value = -46595.5

$ws.Range("H31").Value = 4660.3125
$ws.Range("I31").Value = 4019.5625
$ws.Range("J31").Value = 5301.0625
$ws.Range("K31").Value = 4019.5625
$ws.Range("L31").Value = 5301.0625
$ws.Range("M31").Value = -3724.5625
$ws.Range("N31").Value = -5891.0625

$ws.Range("H34").Value = 4660.3125
$ws.Range("I34").Value = 4019.5625
$ws.Range("J34").Value = 5301.0625
$ws.Range("K34").Value = 4019.5625
$ws.Range("L34").Value = 5301.0625
$ws.Range("M34").Value = -3817.5625
$ws.Range("N34").Value = -5705.0625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2003.8506
$ws.Range("I68").Value = 3213.9429
$ws.Range("J68").Value = 1189.3654
$ws.Range("K68").Value = 9641.8287
$ws.Range("L68").Value = 3568.0962
$ws.Range("M68").Value = -8830.8287
$ws.Range("N68").Value = -5190.0962

$ws.Range("H71").Value = 2003.8506
$ws.Range("I71").Value = 3213.9429
$ws.Range("J71").Value = 1189.3654
$ws.Range("K71").Value = 28925.4861
$ws.Range("L71").Value = 10704.2886
$ws.Range("M71").Value = -24869.4861
$ws.Range("N71").Value = -18816.2886

$ws.Range("H132").Value = 1184.9459
$ws.Range("I132").Value = 896.6667
$ws.Range("J132").Value = 1563.3125
$ws.Range("K132").Value = 8070.0003
$ws.Range("L132").Value = 14069.8125
$ws.Range("M132").Value = -5540.0003
$ws.Range("N132").Value = -19129.8125

$ws.Range("H139").Value = 2866.5186
$ws.Range("I139").Value = 1341
$ws.Range("J139").Value = 3763.8823
$ws.Range("K139").Value = 4023
$ws.Range("L139").Value = 11291.6469
$ws.Range("M139").Value = 1117
$ws.Range("N139").Value = -21571.6469

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3000
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -4996

$ws.Range("H83").Value = 3000
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -24984

$ws.Range("H132").Value = 39283.355
$ws.Range("I132").Value = 2723.5454
$ws.Range("J132").Value = 62939.707
$ws.Range("K132").Value = 8170.6362
$ws.Range("L132").Value = 188819.121
$ws.Range("M132").Value = -5640.6362
$ws.Range("N132").Value = -193879.121

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 545.75
$ws.Range("I22").Value = 471.69232
$ws.Range("J22").Value = 866.6667
$ws.Range("K22").Value = 471.69232
$ws.Range("L22").Value = 866.6667
$ws.Range("M22").Value = -176.69232
$ws.Range("N22").Value = -1456.6667

$ws.Range("H27").Value = 545.75
$ws.Range("I27").Value = 471.69232
$ws.Range("J27").Value = 866.6667
$ws.Range("K27").Value = 471.69232
$ws.Range("L27").Value = 866.6667
$ws.Range("M27").Value = -364.69232
$ws.Range("N27").Value = -1080.6667

$ws.Range("H55").Value = 316.88095
$ws.Range("I55").Value = 238.08333
$ws.Range("J55").Value = 421.94446
$ws.Range("K55").Value = 238.08333
$ws.Range("L55").Value = 421.94446
$ws.Range("M55").Value = -65.08332999999999
$ws.Range("N55").Value = -767.9444599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3106
$ws.Range("I81").Value = 1950
$ws.Range("J81").Value = 3461.6924
$ws.Range("K81").Value = 3900
$ws.Range("L81").Value = 6923.3848
$ws.Range("M81").Value = -2839
$ws.Range("N81").Value = -9045.3848

$ws.Range("H84").Value = 3106
$ws.Range("I84").Value = 1950
$ws.Range("J84").Value = 3461.6924
$ws.Range("K84").Value = 19500
$ws.Range("L84").Value = 34616.924
$ws.Range("M84").Value = -14196
$ws.Range("N84").Value = -45224.924

Write-Output "Applied Tiamat_Profits market-price updates across 8 sheets (27 rows)."
